# Cavezzo.xlsx update: one extra day (2021-02-08 / serial 44235) is inserted
# into the daily series, and the trailing window is extended by one more
# day (2021-03-02 / serial 44257). The 7-day rolling columns (C and D) are
# recomputed accordingly for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 93, pushing old rows 93..113 down to 94..114.
$ws.Rows.Item(93).Insert()

# The freshly inserted row 93 doesn't carry the date-column formatting;
# clone it from the (untouched) row above so A93 keeps the date display.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)

# 2) Write the new row's data (2021-02-08).
$ws.Range("A93").Value2 = 44235
$ws.Range("B93").Value2 = 1
$ws.Range("C93").Value2 = 13
$ws.Range("D93").Value2 = 184.9217638691323

# 3) Row 92 keeps its own date/count, but its rolling totals change now
#    that the new day sits inside later windows.
$ws.Range("C92").Value2 = 12
$ws.Range("D92").Value2 = 170.697012802276

# 4) Recompute the 7-day rolling columns for every shifted row whose
#    window now includes the inserted day (rows 94..112).
$ws.Range("C94").Value2 = 13
$ws.Range("D94").Value2 = 184.9217638691323

$ws.Range("C95").Value2 = 10
$ws.Range("D95").Value2 = 142.2475106685633

$ws.Range("C96").Value2 = 9
$ws.Range("D96").Value2 = 128.022759601707

$ws.Range("C97").Value2 = 9
$ws.Range("D97").Value2 = 128.022759601707

$ws.Range("C98").Value2 = 8
$ws.Range("D98").Value2 = 113.7980085348507

$ws.Range("C99").Value2 = 7
$ws.Range("D99").Value2 = 99.5732574679943

$ws.Range("C100").Value2 = 5
$ws.Range("D100").Value2 = 71.12375533428165

$ws.Range("C101").Value2 = 3
$ws.Range("D101").Value2 = 42.67425320056899

$ws.Range("C102").Value2 = 3
$ws.Range("D102").Value2 = 42.67425320056899

$ws.Range("C103").Value2 = 3
$ws.Range("D103").Value2 = 42.67425320056899

$ws.Range("C104").Value2 = 3
$ws.Range("D104").Value2 = 42.67425320056899

$ws.Range("C105").Value2 = 3
$ws.Range("D105").Value2 = 42.67425320056899

$ws.Range("C106").Value2 = 4
$ws.Range("D106").Value2 = 56.89900426742533

$ws.Range("C107").Value2 = 5
$ws.Range("D107").Value2 = 71.12375533428165

$ws.Range("C108").Value2 = 9
$ws.Range("D108").Value2 = 128.022759601707

$ws.Range("C109").Value2 = 9
$ws.Range("D109").Value2 = 128.022759601707

$ws.Range("C110").Value2 = 12
$ws.Range("D110").Value2 = 170.697012802276

$ws.Range("C111").Value2 = 12
$ws.Range("D111").Value2 = 170.697012802276

$ws.Range("C112").Value2 = 14
$ws.Range("D112").Value2 = 199.1465149359886

# Rows 113 and 114 (post-shift) keep the blank C/D cells they already had
# as rows 112/113 before the insert, so nothing else to do for them.

# 5) Append the new last day (2021-03-02, serial 44257) with 2 new cases,
#    matching the existing style used by the rest of the date column.
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value2 = 44257
$ws.Range("B115").Value2 = 2
